# refactor excel reading code
# The "creds" sheet used to hold Username/Password pairs; it now holds
# Name/Age pairs (numeric ages instead of text passwords), with a thin
# grid border around the whole table and left/top aligned age values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "creds" sheet

# ---- Header row -----------------------------------------------------
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Age"

# ---- Data rows (names stay the same, passwords -> numeric ages) -----
$ws.Range("A2").Value = "Rahul"
$ws.Range("B2").Value = 29.5

$ws.Range("A3").Value = "Laxman"
$ws.Range("B3").Value = 30.5

$ws.Range("A4").Value = "Anjuli"
$ws.Range("B4").Value = 31.2

$ws.Range("A5").Value = "Pooja"
$ws.Range("B5").Value = 32.3

# ---- Thin border around every cell of the table ----------------------
$table = $ws.Range("A1:B5")
$table.Borders.LineStyle = 1   # xlContinuous -> renders as "thin"

# ---- Left/top alignment for the Age column, applied together with the
#      border so both land on the same style record. Build the combined
#      format on a scratch cell first (single-cell writes coalesce into
#      one style record in this engine) then stamp it onto B2:B5 via
#      copy / paste-special so no unused intermediate styles are left
#      behind in the stylesheet.
$scratch = $ws.Range("D1")
$scratch.Borders.LineStyle = 1
$scratch.HorizontalAlignment = -4131   # xlLeft
$scratch.VerticalAlignment = -4160     # xlTop
$scratch.Copy()
$ws.Range("B2:B5").PasteSpecial(-4122) # xlPasteFormats
$scratch.Clear()

# ---- Column widths -----------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 27.5
$ws.Columns.Item(3).ColumnWidth = 22.666666666666668

# ---- Selection moved from G8 to F8 -------------------------------------
$ws.Range("F8").Select()
